# Python_Amphiboles_Test.xlsx -- apply commit "Added element to oxide, Keith Si
# equation, Before messed with labelling":
#   1. Insert a new "Info" sheet as the first sheet explaining where
#      Putirka_Benchmarks pulls its data from.
#   2. On Putirka_Benchmarks, add two new calculated-value columns:
#      SiO2_Calc (BB) and T_C_SiO2_Calc (BC).
#   3. Leave the tab selection on Putirka_Benchmarks (scrolled over to the
#      new columns) instead of Sheet1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Info" sheet, first tab, single explanatory cell.
# ---------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "Info"
$infoSheet.Range("A1").Value = "Putirka_Benchmarks reads from Putirka_Amphibole P-T_v.6_downMay2021.xlsx, the most up-to-date version"
$infoSheet.Move($wb.Worksheets.Item(1))

# ---------------------------------------------------------------------
# 2. New SiO2_Calc / T_C_SiO2_Calc columns on Putirka_Benchmarks.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Putirka_Benchmarks")

$ws.Range("BB1").Value = "SiO2_Calc"
$ws.Range("BC1").Value = "T_C_SiO2_Calc"

$ws.Range("BB2").Value = 57.936794676762176
$ws.Range("BC2").Value = 955.93956668022042
$ws.Range("BB3").Value = 63.369069414151781
$ws.Range("BC3").Value = 911.89573022935838
$ws.Range("BB4").Value = 55.154572604355693
$ws.Range("BC4").Value = 978.97328086757784
$ws.Range("BB5").Value = 62.682265934371905
$ws.Range("BC5").Value = 935.43054369867843
$ws.Range("BB6").Value = 64.774797714210138
$ws.Range("BC6").Value = 914.6345629403919
$ws.Range("BB7").Value = 60.056929592101881
$ws.Range("BC7").Value = 958.5507702281169
$ws.Range("BB8").Value = 62.538583166103997
$ws.Range("BC8").Value = 932.88405708448386
$ws.Range("BB9").Value = 64.0207866165794
$ws.Range("BC9").Value = 908.93954404667261
$ws.Range("BB10").Value = 60.044381225346953
$ws.Range("BC10").Value = 954.78564773000312
$ws.Range("BB11").Value = 53.105290465119687
$ws.Range("BC11").Value = 975.82496756237549
$ws.Range("BB12").Value = 54.037494320022375
$ws.Range("BC12").Value = 964.69486450788281
$ws.Range("BB13").Value = 52.079901953608392
$ws.Range("BC13").Value = 985.92565312692466
$ws.Range("BB14").Value = 57.936794676762176
$ws.Range("BC14").Value = 955.93956668022042
$ws.Range("BB15").Value = 63.369069414151781
$ws.Range("BC15").Value = 911.89573022935838
$ws.Range("BB16").Value = 55.154572604355693
$ws.Range("BC16").Value = 978.97328086757784
$ws.Range("BB17").Value = 72.270201280929896
$ws.Range("BC17").Value = 807.47189104003451
$ws.Range("BB18").Value = 73.502618363842714
$ws.Range("BC18").Value = 793.17204982520741
$ws.Range("BB19").Value = 73.502618363842714
$ws.Range("BC19").Value = 793.17204982520741
$ws.Range("BB20").Value = 65.139424765079326
$ws.Range("BC20").Value = 684.6109011016058
$ws.Range("BB21").Value = 65.525355917219827
$ws.Range("BC21").Value = 690.99385346086331
$ws.Range("BB22").Value = 68.398915624742628
$ws.Range("BC22").Value = 791.56115992271941
$ws.Range("BB23").Value = 63.258429461137496
$ws.Range("BC23").Value = 809.75166252155498

# ---------------------------------------------------------------------
# 3. View state: Putirka_Benchmarks becomes the active/selected tab,
#    scrolled/zoomed over the new columns; Sheet1 no longer selected.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 60
$ws.Range("AK20").Select()
